$wb = $excel.ActiveWorkbook

$wsMAB = $wb.Worksheets.Item("u_MAB")
$wsMAB.Range("B15").Value = 0
$wsMAB.Range("A16").Value = 0
$wsMAB.Range("B27").Value = 1.286967141090104
$wsMAB.Range("A40").Value = 0
$wsMAB.Range("B40").Value = 0
$wsMAB.Range("A49").Value = 0.4691391785999066
$wsMAB.Range("A51").Value = 1.3222106800015
$wsMAB.Range("B51").Value = 0.1180187419164171
$wsMAB.Range("A52").Value = 0.05182702263477304
$wsMAB.Range("B61").Value = 0

$wsEOH = $wb.Worksheets.Item("u_EOH")
$wsEOH.Range("A2").Value = -0.3232560766459875
$wsEOH.Range("A3").Value = -0.5837314435043842

$wsVL = $wb.Worksheets.Item("v_l")
$wsVL.Range("A2").Value = 3319108.402887601
$wsVL.Range("A3").Value = 2728602.857971512
$wsVL.Range("A4").Value = 8442122.906447072
